$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 474 (new rows 475 & 476), shifting the
# existing rows 475:600 down to 477:602 (Excel default shift = down).
$ws.Rows("475:476").Insert()

# Row 475 (new): Packham's Triumph / Especial
$ws.Range("A475").Value = 11
$ws.Range("B475").Value = "Vega Monumental Concepción"
$ws.Range("C475").Value = "Bíobío"
$ws.Range("D475").Value = 44995
$ws.Range("E475").Value = 8
$ws.Range("F475").Value = "Fruta"
$ws.Range("G475").Value = 100104
$ws.Range("H475").Value = "Frutos de pepita"
$ws.Range("I475").Value = 100104005
$ws.Range("J475").Value = "Pera"
$ws.Range("K475").Value = "Packham's Triumph"
$ws.Range("L475").Value = "Especial"
$ws.Range("M475").Value = 250
$ws.Range("N475").Value = 13000
$ws.Range("O475").Value = 13000
$ws.Range("P475").Value = 13000
$ws.Range("Q475").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R475").Value = "Región de O'Higgins"
$ws.Range("S475").Value = 812
$ws.Range("T475").Value = 16

# Row 476 (new): Packham's Triumph / Primera
$ws.Range("A476").Value = 11
$ws.Range("B476").Value = "Vega Monumental Concepción"
$ws.Range("C476").Value = "Bíobío"
$ws.Range("D476").Value = 44995
$ws.Range("E476").Value = 8
$ws.Range("F476").Value = "Fruta"
$ws.Range("G476").Value = 100104
$ws.Range("H476").Value = "Frutos de pepita"
$ws.Range("I476").Value = 100104005
$ws.Range("J476").Value = "Pera"
$ws.Range("K476").Value = "Packham's Triumph"
$ws.Range("L476").Value = "Primera"
$ws.Range("M476").Value = 300
$ws.Range("N476").Value = 12000
$ws.Range("O476").Value = 12000
$ws.Range("P476").Value = 12000
$ws.Range("Q476").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R476").Value = "Región de O'Higgins"
$ws.Range("S476").Value = 750
$ws.Range("T476").Value = 16
